# Updated cryptos list: refresh Price/Volume(1h) figures and re-rank two coin pairs
# (EthereumClassic/InternetComputer swap rows 26-27, ThetaToken/EOS/Cronos shift rows 49-51)
# to mirror the upstream GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", "30.283.39"),
    @("E2", "  +1.99%  "),
    @("D3", "2.093.50"),
    @("E3", "  -0.20%  "),
    @("E4", "  -0.62%  "),
    @("D5", "341.33"),
    @("E5", "  -0.58%  "),
    @("E6", "  -0.63%  "),
    @("D7", "0.5303"),
    @("E7", "  +2.16%  "),
    @("D8", "0.4378"),
    @("E8", "  +0.03%  "),
    @("D9", "54.39"),
    @("E9", "  +1.32%  "),
    @("D10", "0.09347"),
    @("E10", "  +1.40%  "),
    @("D11", "1.174"),
    @("E11", "  +0.63%  "),
    @("D12", "24.66"),
    @("E12", "  +0.19%  "),
    @("D13", "8.574"),
    @("E13", "  +5.08%  "),
    @("D14", "6.876"),
    @("E14", "  +1.37%  "),
    @("D15", "2.044.47"),
    @("E15", "  -1.11%  "),
    @("D16", "101.33"),
    @("E16", "  -1.89%  "),
    @("D17", "0.00001155"),
    @("E17", "  +0.25%  "),
    @("D18", "1.003"),
    @("E18", "  -0.58%  "),
    @("D19", "21.09"),
    @("E19", "  +0.56%  "),
    @("D20", "0.06711"),
    @("E20", "  +0.60%  "),
    @("D21", "6.331"),
    @("E21", "  +1.93%  "),
    @("D22", "1.002"),
    @("E22", "  -0.57%  "),
    @("D23", "30.270.28"),
    @("E23", "  +1.85%  "),
    @("D24", "12.45"),
    @("E24", "  -0.69%  "),
    @("D25", "2.314"),
    @("E25", "  +0.38%  "),
    @("B26", "EthereumClassic"),
    @("C26", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"),
    @("D26", "21.78"),
    @("E26", "  -0.48%  "),
    @("B27", "InternetComputer(DFINITY)"),
    @("C27", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"),
    @("D27", "6.922"),
    @("E27", "  +9.09%  "),
    @("D28", "162.17"),
    @("E28", "  +0.12%  "),
    @("D29", "2.503"),
    @("E29", "  +0.52%  "),
    @("D30", "133.57"),
    @("E30", "  +0.10%  "),
    @("D31", "1.130"),
    @("E31", "  +0.27%  "),
    @("E32", "  +0.18%  "),
    @("E33", "  -1.97%  "),
    @("D34", "6.243"),
    @("E34", "  +0.84%  "),
    @("D35", "3.912"),
    @("E35", "  -1.01%  "),
    @("D36", "10.04"),
    @("E36", "  -3.56%  "),
    @("D37", "0.02618"),
    @("E37", "  +1.67%  "),
    @("D38", "0.06758"),
    @("E38", "  +0.57%  "),
    @("D39", "12.57"),
    @("E39", "  +0.68%  "),
    @("D40", "0.6948"),
    @("E40", "  -0.50%  "),
    @("D41", "1.342"),
    @("E41", "  +1.21%  "),
    @("D42", "0.2213"),
    @("E42", "  -0.01%  "),
    @("D43", "0.6770"),
    @("E43", "  -0.20%  "),
    @("D44", "2.350"),
    @("E44", "  +1.00%  "),
    @("D45", "14.21"),
    @("E45", "  -0.57%  "),
    @("E46", "  -0.52%  "),
    @("D47", "1.274"),
    @("E47", "  +6.11%  "),
    @("D48", "3.631"),
    @("E48", "  +0.20%  "),
    @("B49", "ThetaToken"),
    @("C49", "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"),
    @("D49", "1.207"),
    @("E49", "  +4.30%  "),
    @("B50", "EOS"),
    @("C50", "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"),
    @("D50", "1.212"),
    @("E50", "  -0.34%  "),
    @("B51", "Cronos"),
    @("C51", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"),
    @("D51", "0.07277"),
    @("E51", "  +3.38%  ")
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $cellVal = $u[1]
    $rng = $ws.Range($cellRef)
    # Source workbook stores these as inline text (e.g. "30.283.39", "1.130"),
    # so force Text format first or Excel would coerce numeric-looking strings to numbers.
    $rng.NumberFormat = "@"
    $rng.Value = $cellVal
}
